$wb = $excel.ActiveWorkbook

# --- "Repayment schedule" sheet: insert a new column N (Variable Instalments) ---
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new blank column before the old "Late" column (was N, becomes O).
$ws.Columns.Item(14).EntireColumn.Insert()

# Match the new column's width to column M ("In Advance").
$ws.Columns.Item(14).ColumnWidth = 9.8

# Leave the selection where the author left it and make this the active sheet/tab.
$ws.Range("K15").Select() | Out-Null
